$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2 currently holds "1h 55m" - update it to "end"
$ws.Range("C2").Value = "end"

# Move the active selection from D3 to C3
$ws.Range("C3").Select()
